# Partial resolution of plotting
# The "Aperture of quad" (Aperture / Circular / Radius / 0.005 / m) row that
# already existed once after the focusing-quad (Fquad) parameters and once
# after the defocusing-quad (Dquad) parameters is duplicated so that every
# quad parameter (Fquad Length, Fquad Strength, Dquad Length, Dquad Strength)
# is immediately followed by its own "Aperture of quad" row.
#
# Concretely two new rows are inserted in the "Capture" section:
#   1. Before "Length of focusing quad"  (originally row 15)
#   2. Before "Length of defocusing quad" (originally row 19, row 20 after
#      the first insertion has shifted things down by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insertion #1: new "Aperture of quad" row before the Fquad Length row ---
$ws.Rows("15:15").Insert()
$ws.Range("A18:H18").Copy()              # row 18 is still the original "Aperture of quad" row (unshifted by this insert)
$ws.Range("A15:H15").PasteSpecial(-4122) # xlPasteFormats - carry over the formatting/style only
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Capture"
$ws.Range("C15").Value = "Aperture"
$ws.Range("D15").Value = "Circular"
$ws.Range("E15").Value = "Radius"
$ws.Range("F15").Value = 0.005
$ws.Range("G15").Value = "m"
$ws.Range("H15").Value = "Aperture of quad"

# --- Insertion #2: new "Aperture of quad" row before the Dquad Length row ---
$ws.Rows("20:20").Insert()
$ws.Range("A18:H18").Copy()              # row 18 still holds the original "Aperture of quad" row/format
$ws.Range("A20:H20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Capture"
$ws.Range("C20").Value = "Aperture"
$ws.Range("D20").Value = "Circular"
$ws.Range("E20").Value = "Radius"
$ws.Range("F20").Value = 0.005
$ws.Range("G20").Value = "m"
$ws.Range("H20").Value = "Aperture of quad"

$ws.Application.CutCopyMode = $false

# Match the saved selection recorded in the authored workbook.
[void]$ws.Range("A20:H20").Select()
